$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix typo in A15: "הילס שולויס" -> "הילה שולויס"
$ws.Range("A15").Value = "הילה שולויס"

# Rename "ליאם מלכה" -> "ליאם דיין " (with trailing space) in A52 and A58
$ws.Range("A52").Value = "ליאם דיין "
$ws.Range("A58").Value = "ליאם דיין "

# Update the active selection to A10
$ws.Range("A10").Select()
